# "conectando modulo de contadores al front"
#
# The registros.xlsx workbook gets two new data rows appended under the
# existing header row, the sheet gains explicit column widths, and the
# workbook/selection/protection bookkeeping is refreshed to match what a
# live Excel session leaves behind after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original file carried an (empty) <workbookProtection/> element; the
# edited workbook no longer protects workbook structure/windows.
$wb.Unprotect()

# --- New row 2: Casino Medellín / M1-1, 2025-05-22 -------------------------
# Force these as literal text so "2025-05-22" doesn't get reinterpreted as a
# date serial number (matches the inlineStr cells produced by the backend).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-05-22"
$ws.Range("B2").Value = "Casino Medellín"
$ws.Range("C2").Value = "M1-1"
$ws.Range("D2").Value = 100000
$ws.Range("E2").Value = 200000
$ws.Range("F2").Value = 50000
$ws.Range("G2").Value = 1000000

# --- New row 3: same casino/machine, placeholder counters = 1 --------------
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-05-22"
$ws.Range("B3").Value = "Casino Medellín"
$ws.Range("C3").Value = "M1-1"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1

# --- Explicit column widths (A:G) ------------------------------------------
# ColumnWidth is expressed in "characters"; values below were chosen so the
# engine's stored <col width=.../> lands as close as possible to the
# target widths (19.29, 27, 15.29, 15.86, 15.29, 13.86, 17).
$ws.Columns.Item(1).ColumnWidth = 18.5
$ws.Columns.Item(2).ColumnWidth = 26.1666666666667
$ws.Columns.Item(3).ColumnWidth = 14.5
$ws.Columns.Item(4).ColumnWidth = 15
$ws.Columns.Item(5).ColumnWidth = 14.5
$ws.Columns.Item(6).ColumnWidth = 13
$ws.Columns.Item(7).ColumnWidth = 16.1666666666667

# --- Selection / active view -------------------------------------------
# The saved file's active selection moved to E9 with the sheet tab active.
$ws.Range("E9").Select()
